# rerun corona results with larger ds
# Rebuild the whole results sheet (anchor-word confidence tables) with the
# refreshed numbers produced by the larger dataset run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - all old content/formatting is replaced.
$ws.Cells.Clear()

# ---- "negative" anchor table (columns A:H) ----
$ws.Range("A1").Value = "negative"

$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"

$tableA = @(
    @('crude', 0.8823529411764706, 30, 30, 0, 1, $false, 4),
    @('crisis', 0.6335616438356164, 185, 185, 0, 1, $false, 107),
    @('panic', 0.2209302325581395, 114, 114, 0, 1, $false, 402),
    @('sc', 0.1851851851851852, 35, 35, 0, 1, $false, 154),
    @('no', 0.09722222222222222, 35, 35, 0, 1, $false, 325)
)

$rowIdx = 3
foreach ($row in $tableA) {
    $ws.Cells.Item($rowIdx, 1).Value = $row[0]
    $ws.Cells.Item($rowIdx, 2).Value = $row[1]
    $ws.Cells.Item($rowIdx, 3).Value = $row[2]
    $ws.Cells.Item($rowIdx, 4).Value = $row[3]
    $ws.Cells.Item($rowIdx, 5).Value = $row[4]
    $ws.Cells.Item($rowIdx, 6).Value = $row[5]
    $ws.Cells.Item($rowIdx, 7).Value = $row[6]
    $ws.Cells.Item($rowIdx, 8).Value = $row[7]
    $rowIdx++
}

# ---- "positive" anchor table (columns J:Q) ----
$ws.Range("J1").Value = "positive"

$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"

$tableJ = @(
    @('best', 0.9491525423728814, 56, 56, 1, 0, $false, 3),
    @('interesting', 0.9393939393939394, 31, 31, 1, 0, $false, 2),
    @('love', 0.9130434782608695, 42, 42, 1, 0, $false, 4),
    @('great', 0.8482142857142857, 95, 95, 1, 0, $false, 17),
    @('won', 0.8461538461538461, 33, 33, 1, 0, $false, 6),
    @('thanks', 0.7926829268292683, 65, 65, 1, 0, $false, 17),
    @('free', 0.7583333333333333, 91, 91, 1, 0, $false, 29),
    @('thank', 0.7578125, 97, 97, 1, 0, $false, 31),
    @('positive', 0.7413793103448276, 43, 43, 1, 0, $false, 15),
    @('special', 0.7222222222222222, 26, 26, 1, 0, $false, 10),
    @('safe', 0.7183098591549296, 102, 102, 1, 0, $false, 40),
    @('support', 0.7169811320754716, 76, 76, 1, 0, $false, 30),
    @('good', 0.7125, 114, 114, 1, 0, $false, 46),
    @('safety', 0.6470588235294118, 33, 33, 1, 0, $false, 18),
    @('well', 0.6276595744680851, 59, 59, 1, 0, $false, 35),
    @('better', 0.5873015873015873, 37, 37, 1, 0, $false, 26),
    @('fresh', 0.5833333333333334, 28, 28, 1, 0, $false, 20),
    @('relief', 0.58, 29, 29, 1, 0, $false, 21),
    @('heroes', 0.5319148936170213, 25, 25, 1, 0, $false, 22),
    @('hand', 0.5195822454308094, 199, 199, 1, 0, $false, 184),
    @('care', 0.4943820224719101, 44, 44, 1, 0, $false, 45),
    @('help', 0.4779661016949153, 141, 141, 1, 0, $false, 154),
    @('like', 0.4647058823529412, 158, 158, 1, 0, $false, 182),
    @('increase', 0.3717948717948718, 29, 29, 1, 0, $false, 49),
    @('please', 0.3472803347280335, 83, 83, 1, 0, $false, 156),
    @('grocery', 0.03551609322974473, 32, 32, 1, 0, $false, 869),
    @('store', 0.03355704697986577, 30, 30, 1, 0, $false, 864),
    @('you', 0.02666666666666667, 32, 32, 1, 0, $false, 1168),
    @(',', 0.01884473576403113, 46, 50, 0.92, 0.07999999999999996, $true, 2395),
    @('for', 0.01756587202007528, 28, 28, 1, 0, $false, 1566),
    @('and', 0.0157244477723699, 42, 44, 0.95, 0.05000000000000004, $true, 2629),
    @('to', 0.01132947976878613, 49, 52, 0.94, 0.06000000000000005, $true, 4276),
    @('a', 0.01127004768097096, 26, 28, 0.93, 0.06999999999999995, $true, 2281),
    @('.', 0.009617311160088158, 48, 53, 0.91, 0.08999999999999997, $true, 4943),
    @('co', 0.009017713365539453, 28, 30, 0.93, 0.06999999999999995, $true, 3077),
    @('the', 0.008137957760124007, 42, 46, 0.91, 0.08999999999999997, $true, 5119)
)

$rowIdx = 3
foreach ($row in $tableJ) {
    $ws.Cells.Item($rowIdx, 10).Value = $row[0]
    $ws.Cells.Item($rowIdx, 11).Value = $row[1]
    $ws.Cells.Item($rowIdx, 12).Value = $row[2]
    $ws.Cells.Item($rowIdx, 13).Value = $row[3]
    $ws.Cells.Item($rowIdx, 14).Value = $row[4]
    $ws.Cells.Item($rowIdx, 15).Value = $row[5]
    $ws.Cells.Item($rowIdx, 16).Value = $row[6]
    $ws.Cells.Item($rowIdx, 17).Value = $row[7]
    $rowIdx++
}

# ---- Formatting: bold + thin border + centered/top-aligned header style ----
# Applies to the header rows (A2:H2 and J2:Q2) and to the "name" column of
# every data row in both tables (A3:A7 and J3:J38).
$headerStyleRanges = @("A2:H2", "J2:Q2", "A3:A7", "J3:J38")
foreach ($addr in $headerStyleRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}
